# Adds a "cfop" breakdown to the workbook:
#   1. A new "cfop" column (G) on the "PI hours" sheet, mirroring the
#      existing "app" column's layout/formatting.
#   2. A new "cfop hours" worksheet (after "unit(accumulative) hours"),
#      built the same way as "department hours" / "unit(accumulative) hours".

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. "PI hours" sheet -> add the "cfop" column
# ---------------------------------------------------------------------
$piSheet = $wb.Worksheets.Item("PI hours")

# Copy header formatting from the neighboring "app" header cell (F1) to
# the new "cfop" header cell (G1).
$piSheet.Range("F1").Copy()
$piSheet.Range("G1").PasteSpecial(-4122)   # xlPasteFormats

$piSheet.Cells.Item(1, 7).Value = "cfop"

$piSheet.Cells.Item(2, 7).Value = "['cfop_GC']"
$piSheet.Cells.Item(3, 7).Value = "['cfop_HUTCHINSON']"
$piSheet.Cells.Item(4, 7).Value = "['cfop_NH']"
$piSheet.Cells.Item(5, 7).Value = "['cfop_CHOUDHURY', 'cfop_RRC']"

# ---------------------------------------------------------------------
# 2. New "cfop hours" worksheet, placed after "unit(accumulative) hours"
# ---------------------------------------------------------------------
$unitSheet = $wb.Worksheets.Item("unit(accumulative) hours")
$cfopSheet = $wb.Worksheets.Add($null, $unitSheet)
$cfopSheet.Name = "cfop hours"

# Copy the header / id-column formatting from the "unit(accumulative)
# hours" sheet so the new sheet matches the look of the others.
$unitSheet.Range("B1:D1").Copy()
$cfopSheet.Range("B1:D1").PasteSpecial(-4122)   # xlPasteFormats

$unitSheet.Range("A2:A6").Copy()
$cfopSheet.Range("A2:A6").PasteSpecial(-4122)   # xlPasteFormats

$cfopSheet.Cells.Item(1, 2).Value = "cfop"
$cfopSheet.Cells.Item(1, 3).Value = "hours"
$cfopSheet.Cells.Item(1, 4).Value = "percentage"

$cfopSheet.Cells.Item(2, 1).Value = 0
$cfopSheet.Cells.Item(2, 2).Value = "cfop_GC"
$cfopSheet.Cells.Item(2, 3).Value = 95
$cfopSheet.Cells.Item(2, 4).Value = 47.73869346733668

$cfopSheet.Cells.Item(3, 1).Value = 1
$cfopSheet.Cells.Item(3, 2).Value = "cfop_HUTCHINSON"
$cfopSheet.Cells.Item(3, 3).Value = 78
$cfopSheet.Cells.Item(3, 4).Value = 39.19597989949749

$cfopSheet.Cells.Item(4, 1).Value = 2
$cfopSheet.Cells.Item(4, 2).Value = "cfop_NH"
$cfopSheet.Cells.Item(4, 3).Value = 14
$cfopSheet.Cells.Item(4, 4).Value = 7.035175879396985

$cfopSheet.Cells.Item(5, 1).Value = 3
$cfopSheet.Cells.Item(5, 2).Value = "cfop_RRC"
$cfopSheet.Cells.Item(5, 3).Value = 8
$cfopSheet.Cells.Item(5, 4).Value = 4.020100502512562

$cfopSheet.Cells.Item(6, 1).Value = 4
$cfopSheet.Cells.Item(6, 2).Value = "cfop_CHOUDHURY"
$cfopSheet.Cells.Item(6, 3).Value = 4
$cfopSheet.Cells.Item(6, 4).Value = 2.010050251256281

$piSheet.Select()
